# edit.ps1 -- applies the "added slight changes to Proposal and Data Analysis"
# commit to Data_Analysis.docx via the Word COM object model.

$d = $word.ActiveDocument

# wdReplaceAll = 2, wdFindContinue = 1, wdAlignParagraphCenter = 1, wdAlignParagraphLeft = 0

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 1. Simple literal text replacements (run-splitting in the source diff does
#    not change the rendered text, so a plain Find/Replace is equivalent).
# ---------------------------------------------------------------------------

Replace-Text "were neck in neck. The gap" "were neck and neck. The gap"

Replace-Text "is quite alarming and will remain a common theme" `
             "is quite large. This will, unfortunately for the region of San Antonio, remain a common theme"

Replace-Text "the average scores for 2017 (improved" "the average scores for the class of 2017 (improved"

Replace-Text "in the later years that the ACT participation rate grew" `
             "in the later years that ACT participation grew"

Replace-Text ", but it is quite concerning that students are not being encouraged to take the ACT" `
             ", but it is also quite concerning that students are not more encouraged to take the ACT"

Replace-Text "participation rates contain similarly positive correlation" `
             "participation rates contain positive correlation"

Replace-Text "take more exams and potentially earn more college credit." `
             "take more exams and earn more college credit."

Replace-Text "the second most opportunity which is quite impressive" `
             "the second most opportunity, which is quite impressive"

Replace-Text "our next analysis of the Wealth/ADA feature of school districts." `
             "our next analysis of the wealth per average daily attendance (“Wealth/ADA”) feature of school districts."

Replace-Text "there was indeed a positive correlation (.34). The cluster" `
             "there was indeed a positive correlation. The cluster"

Replace-Text "in Fort Worth and Houston. It should be interesting view how Wealth/ADA correlated with college enrollment percentage over the years as well. Let’s take a look. " `
             "in Fort Worth and Houston. Overall, it appears that Wealth/ADA has been increasing in recent years. Something we could later choose to explore is the average property tax for homes in each region. "

Replace-Text "college enrollment percentage over the years as well. Let’s take a look. " `
             "college enrollment percentage for the different class years as well. It would not be unreasonable to hypothesize that a student’s ability to enroll into college increases with wealth as college is not free in the United States. "

Replace-Text "Though it shouldn’t be mistaken for the only factor determining college enrollment percentage, Wealth/ADA certainly contains a positive correlation with college enrollment percentage (.45). As discussed before, this can partially be attributed to “wealthier” students having more access to college even if they do not always have the best scores and aren’t necessarily as prepared to handle the rigors of college. " `
             "Though it certainly isn’t the only factor influencing college enrollment percentage, Wealth/ADA does indeed contain a positive correlation with college enrollment percentage. As mentioned before, this can be attributed to “wealthier” students having more access to college. "

Replace-Text "taking advantage of quality education in testing" `
             "taking advantage of the quality education made available to them in testing"

Replace-Text "it’s not that ridiculous to assume that many of the parents" `
             "it’s not egregious to assume that many of the parents"

Replace-Text "within four years as we could expect many of these unprepared students to fail out of college. " `
             "within four years. "

Replace-Text "The students attending college from less fortunate areas will most likely be intelligent ones who earned scholarship and will have greater chance of earning their degree (more prepared), hurting the positive correlation as well. Let’s take a look at the real historical data from 2011 – 2014 (college graduation year: 2015 – 2018) I collected as a comparison. " `
             "For the less fortunate areas, the students who attended college most likely earned scholarship. We could assume these students contained a greater chance of earning their degree as they were more prepared, hurting the positive correlation as well. To test our assumptions, let’s take a look at the actual data from the classes of 2011 – 2014 (college graduation year: 2015 – 2018)."

Replace-Text "one would expect assuming the hypothetical situation described in the exercise. These wealthier areas tend to attract well-educated families and teachers, so it’s not all that surprising to see that students from these areas tend to do a well in earning their degree percentage-wise. " `
             "one would expect from the assumptions made in the exercise above. These wealthier areas tend to attract well-educated families and teachers, so it’s not all that surprising to see that these school districts contained a decent percentage of their students go on to earn a college degree. "

Replace-Text "(Wealth/ADA Affect on College Graduation %)" "(Wealth/ADA Effect on College Graduation %)"

Write-Output "Phase 1 (text replacements) complete"
